# Update gh-pages to output generated at 456a3b4
# Refresh scraped "want to go" counts (column F) across the sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 15880
$ws1.Range("F8").Value = 712
$ws1.Range("F18").Value = 208
$ws1.Range("F29").Value = 309
$ws1.Range("F30").Value = 35
$ws1.Range("F33").Value = 65

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 74

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 15880
$ws4.Range("F8").Value = 712
$ws4.Range("F18").Value = 208
$ws4.Range("F29").Value = 309
$ws4.Range("F30").Value = 35
$ws4.Range("F32").Value = 74
$ws4.Range("F35").Value = 65
